$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Notes")

# Update the Description line (A2)
$ws.Range("A2").Value = "Description: Urban And Rural Population (%)"

# Update the Source line (A4)
$ws.Range("A4").Value = "Source: National population and Housing census 2014: Provisional Results - Uganda Bureau of Statistics"

# Insert a new row after the Source line for the Source-link
$ws.Range("A5").EntireRow.Insert()
$ws.Range("A5").Value = "Source-link: http://www.ubos.org/onlinefiles/uploads/ubos/NPHC/NPHC%202014%20PROVISIONAL%20RESULTS%20REPORT.pdf"

# Update the license line (now at A12 after the insertion above)
$ws.Range("A12").Value = "It is licensed under a Creative Commons Attribution 4.0 International license."

# Insert a new row after the license line for the licensing info link
$ws.Range("A13").EntireRow.Insert()
$ws.Range("A13").Value = "More information on licensing is available here: https://creativecommons.org/licenses/by/4.0/"
